$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30586
$ws.Range("H13").Value = 8000
$ws.Range("J13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("N13").Value = -8338
$ws.Range("H17").Value = 2107
$ws.Range("I17").Value = 900
$ws.Range("K17").Value = 2700
$ws.Range("M17").Value = -2532
$ws.Range("H19").Value = 898.4286
$ws.Range("I19").Value = 1211.5
$ws.Range("J19").Value = 773.2
$ws.Range("K19").Value = 1211.5
$ws.Range("L19").Value = 773.2
$ws.Range("M19").Value = -1036.5
$ws.Range("N19").Value = -1123.2
$ws.Range("H28").Value = 32459.312
$ws.Range("J28").Value = 882.8333
$ws.Range("L28").Value = 882.8333
$ws.Range("N28").Value = -1852.8333
$ws.Range("H34").Value = 30000
$ws.Range("J34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30406
$ws.Range("H36").Value = 30000
$ws.Range("J36").Value = 30000
$ws.Range("L36").Value = 30000
$ws.Range("N36").Value = -31430
$ws.Range("H43").Value = 1706.4667
$ws.Range("I43").Value = 1733.3334
$ws.Range("K43").Value = 1733.3334
$ws.Range("M43").Value = -1664.3334
$ws.Range("H54").Value = 39084
$ws.Range("J54").Value = 39084
$ws.Range("L54").Value = 39084
$ws.Range("N54").Value = -40056
$ws.Range("H55").Value = 194.75
$ws.Range("I55").Value = 10.5
$ws.Range("K55").Value = 10.5
$ws.Range("M55").Value = 203.5
$ws.Range("H74").Value = 10712.167
$ws.Range("I74").Value = 8319.5
$ws.Range("K74").Value = 8319.5
$ws.Range("M74").Value = -7383.5
$ws.Range("H76").Value = 58830364
$ws.Range("I76").Value = 6817.7
$ws.Range("K76").Value = 6817.7
$ws.Range("M76").Value = -6502.7
$ws.Range("H77").Value = 10712.167
$ws.Range("I77").Value = 8319.5
$ws.Range("K77").Value = 41597.5
$ws.Range("M77").Value = -36917.5
$ws.Range("H79").Value = 58830364
$ws.Range("I79").Value = 6817.7
$ws.Range("K79").Value = 6817.7
$ws.Range("M79").Value = -5725.7
$ws.Range("H97").Value = 4368.4287
$ws.Range("J97").Value = 4368.4287
$ws.Range("L97").Value = 13105.2861
$ws.Range("N97").Value = -14097.2861
$ws.Range("H111").Value = 57105.445
$ws.Range("I111").Value = 78168.69500000001
$ws.Range("J111").Value = 2341
$ws.Range("K111").Value = 234506.085
$ws.Range("L111").Value = 7023
$ws.Range("M111").Value = -231439.085
$ws.Range("N111").Value = -13157
$ws.Range("H132").Value = 1600.9333
$ws.Range("I132").Value = 1513.0927
$ws.Range("K132").Value = 4539.2781
$ws.Range("M132").Value = -2009.2781
$ws.Range("H135").Value = 2125.9
$ws.Range("I135").Value = 1917.7778
$ws.Range("K135").Value = 17260.0002
$ws.Range("M135").Value = -14725.0002
$ws.Range("H137").Value = 5231.08
$ws.Range("I137").Value = 4515.2666
$ws.Range("K137").Value = 13545.7998
$ws.Range("M137").Value = -10995.7998
$ws.Range("H138").Value = 4821.8
$ws.Range("I138").Value = 3406.7727
$ws.Range("J138").Value = 5220.91
$ws.Range("K138").Value = 10220.3181
$ws.Range("L138").Value = 15662.73
$ws.Range("M138").Value = -5080.3181
$ws.Range("N138").Value = -25942.73
$ws.Range("H141").Value = 10367.941
$ws.Range("I141").Value = 10265.9375
$ws.Range("K141").Value = 30797.8125
$ws.Range("M141").Value = -25617.8125

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 41630.84
$ws.Range("I2").Value = 72603.71000000001
$ws.Range("K2").Value = 72603.71000000001
$ws.Range("M2").Value = -72490.71000000001
$ws.Range("H32").Value = 8387.41
$ws.Range("I32").Value = 5261.202
$ws.Range("J32").Value = 24800
$ws.Range("K32").Value = 5261.202
$ws.Range("L32").Value = 24800
$ws.Range("M32").Value = -4974.202
$ws.Range("N32").Value = -25374
$ws.Range("H45").Value = 2987
$ws.Range("I45").Value = 1477.25
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1477.25
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1100.25
$ws.Range("N45").Value = -5754
$ws.Range("H61").Value = 4322.6665
$ws.Range("I61").Value = 3904.907
$ws.Range("J61").Value = 7915.4
$ws.Range("K61").Value = 3904.907
$ws.Range("L61").Value = 7915.4
$ws.Range("M61").Value = -3692.907
$ws.Range("N61").Value = -8339.4
$ws.Range("H63").Value = 6368.143
$ws.Range("I63").Value = 3714.8572
$ws.Range("J63").Value = 9021.429
$ws.Range("K63").Value = 3714.8572
$ws.Range("L63").Value = 9021.429
$ws.Range("M63").Value = -3028.8572
$ws.Range("N63").Value = -10393.429
$ws.Range("H66").Value = 6368.143
$ws.Range("I66").Value = 3714.8572
$ws.Range("J66").Value = 9021.429
$ws.Range("K66").Value = 18574.286
$ws.Range("L66").Value = 45107.145
$ws.Range("M66").Value = -15142.286
$ws.Range("N66").Value = -51971.145
$ws.Range("H74").Value = 3127.52
$ws.Range("I74").Value = 2400.3076
$ws.Range("J74").Value = 3915.3333
$ws.Range("K74").Value = 2400.3076
$ws.Range("L74").Value = 3915.3333
$ws.Range("M74").Value = -1526.3076
$ws.Range("N74").Value = -5663.3333
$ws.Range("H77").Value = 3127.52
$ws.Range("I77").Value = 2400.3076
$ws.Range("J77").Value = 3915.3333
$ws.Range("K77").Value = 12001.538
$ws.Range("L77").Value = 19576.6665
$ws.Range("M77").Value = -7633.538
$ws.Range("N77").Value = -28312.6665
$ws.Range("H102").Value = 1694.2941
$ws.Range("I102").Value = 1648.6129
$ws.Range("K102").Value = 1648.6129
$ws.Range("M102").Value = -26.61290000000008
$ws.Range("H116").Value = 41630.84
$ws.Range("I116").Value = 72603.71000000001
$ws.Range("K116").Value = 72603.71000000001
$ws.Range("M116").Value = -70309.71000000001
$ws.Range("H132").Value = 4309.6377
$ws.Range("I132").Value = 3296.76
$ws.Range("J132").Value = 6975.1055
$ws.Range("K132").Value = 9890.280000000001
$ws.Range("L132").Value = 20925.3165
$ws.Range("M132").Value = -7360.280000000001
$ws.Range("N132").Value = -25985.3165
$ws.Range("H136").Value = 4322.6665
$ws.Range("I136").Value = 3904.907
$ws.Range("J136").Value = 7915.4
$ws.Range("K136").Value = 11714.721
$ws.Range("L136").Value = 23746.2
$ws.Range("M136").Value = -9164.721000000001
$ws.Range("N136").Value = -28846.2

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 41630.84
$ws.Range("I3").Value = 72603.71000000001
$ws.Range("K3").Value = 72603.71000000001
$ws.Range("M3").Value = -72489.71000000001
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H86").Value = 1065233
$ws.Range("I86").Value = 2127342.5
$ws.Range("K86").Value = 2127342.5
$ws.Range("M86").Value = -2126219.5
$ws.Range("H89").Value = 1065233
$ws.Range("I89").Value = 2127342.5
$ws.Range("K89").Value = 10636712.5
$ws.Range("M89").Value = -10631096.5
$ws.Range("H105").Value = 40356.31
$ws.Range("I105").Value = 49150.383
$ws.Range("J105").Value = 3421.2
$ws.Range("K105").Value = 49150.383
$ws.Range("L105").Value = 3421.2
$ws.Range("M105").Value = -47403.383
$ws.Range("N105").Value = -6915.2
$ws.Range("H107").Value = 8589
$ws.Range("I107").Value = 10246.728
$ws.Range("J107").Value = 2510.6667
$ws.Range("K107").Value = 10246.728
$ws.Range("L107").Value = 2510.6667
$ws.Range("M107").Value = -8326.727999999999
$ws.Range("N107").Value = -6350.6667
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H134").Value = 20844.72
$ws.Range("I134").Value = 4228.5615
$ws.Range("J134").Value = 257625
$ws.Range("K134").Value = 12685.6845
$ws.Range("L134").Value = 772875
$ws.Range("M134").Value = -10150.6845
$ws.Range("N134").Value = -777945
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H10").Value = 1750
$ws.Range("J10").Value = 1750
$ws.Range("L10").Value = 1750
$ws.Range("N10").Value = -2028
$ws.Range("H13").Value = 4899.6665
$ws.Range("J13").Value = 4899.6665
$ws.Range("L13").Value = 4899.6665
$ws.Range("N13").Value = -5177.6665
$ws.Range("H16").Value = 14879.125
$ws.Range("I16").Value = 5576
$ws.Range("J16").Value = 80001
$ws.Range("K16").Value = 5576
$ws.Range("L16").Value = 80001
$ws.Range("M16").Value = -5289
$ws.Range("N16").Value = -80575
$ws.Range("H31").Value = 126780.875
$ws.Range("I31").Value = 1569.4
$ws.Range("J31").Value = 335466.66
$ws.Range("K31").Value = 1569.4
$ws.Range("L31").Value = 335466.66
$ws.Range("M31").Value = -1274.4
$ws.Range("N31").Value = -336056.66
$ws.Range("H34").Value = 126780.875
$ws.Range("I34").Value = 1569.4
$ws.Range("J34").Value = 335466.66
$ws.Range("K34").Value = 1569.4
$ws.Range("L34").Value = 335466.66
$ws.Range("M34").Value = -1367.4
$ws.Range("N34").Value = -335870.66
$ws.Range("H50").Value = 26048.309
$ws.Range("I50").Value = 17082.8
$ws.Range("K50").Value = 17082.8
$ws.Range("M50").Value = -16457.8
$ws.Range("H58").Value = 1611
$ws.Range("I58").Value = 1292.3889
$ws.Range("K58").Value = 1292.3889
$ws.Range("M58").Value = -1089.3889
$ws.Range("H99").Value = 5139.1816
$ws.Range("J99").Value = 4587.3335
$ws.Range("L99").Value = 4587.3335
$ws.Range("N99").Value = -7583.3335
$ws.Range("H104").Value = 57968.5
$ws.Range("J104").Value = 57968.5
$ws.Range("L104").Value = 57968.5
$ws.Range("N104").Value = -63210.5
$ws.Range("H107").Value = 916.8946999999999
$ws.Range("I107").Value = 619.25
$ws.Range("J107").Value = 2504.3333
$ws.Range("K107").Value = 619.25
$ws.Range("L107").Value = 2504.3333
$ws.Range("M107").Value = 1300.75
$ws.Range("N107").Value = -6344.3333
$ws.Range("H113").Value = 14879.125
$ws.Range("I113").Value = 5576
$ws.Range("J113").Value = 80001
$ws.Range("K113").Value = 5576
$ws.Range("L113").Value = 80001
$ws.Range("M113").Value = -3406
$ws.Range("N113").Value = -84341
$ws.Range("H126").Value = 5139.1816
$ws.Range("J126").Value = 4587.3335
$ws.Range("L126").Value = 13762.0005
$ws.Range("N126").Value = -18702.0005
$ws.Range("H132").Value = 2610.923
$ws.Range("I132").Value = 2191.1924
$ws.Range("K132").Value = 6573.5772
$ws.Range("M132").Value = -4043.5772
$ws.Range("H136").Value = 1611
$ws.Range("I136").Value = 1292.3889
$ws.Range("K136").Value = 3877.1667
$ws.Range("M136").Value = -1327.1667

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H33").Value = 3086566
$ws.Range("I33").Value = 6172936.5
$ws.Range("K33").Value = 37037619
$ws.Range("M33").Value = -37037336
$ws.Range("H98").Value = 5231.5713
$ws.Range("I98").Value = 4321.5
$ws.Range("J98").Value = 5914.125
$ws.Range("K98").Value = 12964.5
$ws.Range("L98").Value = 17742.375
$ws.Range("M98").Value = -11466.5
$ws.Range("N98").Value = -20738.375
$ws.Range("H113").Value = 1951207.8
$ws.Range("J113").Value = 1995.0555
$ws.Range("L113").Value = 5985.166499999999
$ws.Range("N113").Value = -10325.1665
$ws.Range("H114").Value = 93336.09
$ws.Range("J114").Value = 204210.2
$ws.Range("L114").Value = 612630.6000000001
$ws.Range("N114").Value = -619138.6000000001
$ws.Range("H122").Value = 1406.0769
$ws.Range("I122").Value = 900
$ws.Range("K122").Value = 8100
$ws.Range("M122").Value = -5650
$ws.Range("H131").Value = 68376.94
$ws.Range("I131").Value = 202044
$ws.Range("J131").Value = 42671.73
$ws.Range("K131").Value = 606132
$ws.Range("L131").Value = 128015.19
$ws.Range("M131").Value = -601092
$ws.Range("N131").Value = -138095.19
$ws.Range("H137").Value = 4721.524
$ws.Range("J137").Value = 7399.6665
$ws.Range("L137").Value = 22198.9995
$ws.Range("N137").Value = -32398.9995
$ws.Range("H140").Value = 5012.8213
$ws.Range("I140").Value = 3104.2
$ws.Range("K140").Value = 9312.599999999999
$ws.Range("M140").Value = -4132.599999999999

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H48").Value = 15300
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H75").Value = 30299.5
$ws.Range("J75").Value = 30299.5
$ws.Range("L75").Value = 30299.5
$ws.Range("N75").Value = -32047.5
$ws.Range("H78").Value = 30299.5
$ws.Range("J78").Value = 30299.5
$ws.Range("L78").Value = 90898.5
$ws.Range("N78").Value = -99634.5
$ws.Range("H113").Value = 1115383.5
$ws.Range("J113").Value = 4916.5
$ws.Range("L113").Value = 4916.5
$ws.Range("N113").Value = -9256.5
$ws.Range("H122").Value = 3619
$ws.Range("I122").Value = 3783.4285
$ws.Range("J122").Value = 2468
$ws.Range("K122").Value = 11350.2855
$ws.Range("L122").Value = 7404
$ws.Range("M122").Value = -8900.2855
$ws.Range("N122").Value = -12304
$ws.Range("H132").Value = 28550.834
$ws.Range("I132").Value = 4519.4062
$ws.Range("K132").Value = 13558.2186
$ws.Range("M132").Value = -11028.2186

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 5399.6665
$ws.Range("I7").Value = 3999.6667
$ws.Range("K7").Value = 3999.6667
$ws.Range("M7").Value = -3887.6667
$ws.Range("H22").Value = 1688.6666
$ws.Range("I22").Value = 1688.6666
$ws.Range("K22").Value = 1688.6666
$ws.Range("M22").Value = -1393.6666
$ws.Range("H27").Value = 1688.6666
$ws.Range("I27").Value = 1688.6666
$ws.Range("K27").Value = 1688.6666
$ws.Range("M27").Value = -1581.6666
$ws.Range("H35").Value = 7385.077
$ws.Range("I35").Value = 5878.6
$ws.Range("J35").Value = 12406.667
$ws.Range("K35").Value = 5878.6
$ws.Range("L35").Value = 12406.667
$ws.Range("M35").Value = -5542.6
$ws.Range("N35").Value = -13078.667
$ws.Range("H46").Value = 3050.4167
$ws.Range("I46").Value = 3116.8333
$ws.Range("K46").Value = 3116.8333
$ws.Range("M46").Value = -2928.8333
$ws.Range("H99").Value = 61438
$ws.Range("I99").Value = 50000
$ws.Range("J99").Value = 72876
$ws.Range("K99").Value = 50000
$ws.Range("L99").Value = 72876
$ws.Range("M99").Value = -47005
$ws.Range("N99").Value = -78866
$ws.Range("H122").Value = 7967.357
$ws.Range("I122").Value = 7695
$ws.Range("K122").Value = 23085
$ws.Range("M122").Value = -20635
$ws.Range("H126").Value = 5399.6665
$ws.Range("I126").Value = 3999.6667
$ws.Range("K126").Value = 11999.0001
$ws.Range("M126").Value = -9529.000100000001
$ws.Range("H132").Value = 4241.243
$ws.Range("I132").Value = 4140.5454
$ws.Range("K132").Value = 12421.6362
$ws.Range("M132").Value = -9891.636200000001
$ws.Range("H136").Value = 12266.723
$ws.Range("I136").Value = 9088.444
$ws.Range("J136").Value = 15445
$ws.Range("K136").Value = 27265.332
$ws.Range("L136").Value = 46335
$ws.Range("M136").Value = -24715.332
$ws.Range("N136").Value = -51435
$ws.Range("H138").Value = 93500
$ws.Range("I138").Value = 93500
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 93500
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -88360
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 48999.75
$ws.Range("J141").Value = 48999.75
$ws.Range("L141").Value = 48999.75
$ws.Range("N141").Value = -59359.75

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H4").Value = 15639777
$ws.Range("I4").Value = 18310.666
$ws.Range("J4").Value = 62504176
$ws.Range("K4").Value = 18310.666
$ws.Range("L4").Value = 62504176
$ws.Range("M4").Value = -18197.666
$ws.Range("N4").Value = -62504402
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576
$ws.Range("H18").Value = 29001.084
$ws.Range("J18").Value = 29001.084
$ws.Range("L18").Value = 29001.084
$ws.Range("N18").Value = -29347.084
$ws.Range("H62").Value = 98863.73
$ws.Range("J62").Value = 7812.5
$ws.Range("L62").Value = 7812.5
$ws.Range("N62").Value = -9060.5
$ws.Range("H65").Value = 98863.73
$ws.Range("J65").Value = 7812.5
$ws.Range("L65").Value = 39062.5
$ws.Range("N65").Value = -45302.5
$ws.Range("H107").Value = 1983.1666
$ws.Range("J107").Value = 979.8
$ws.Range("L107").Value = 2939.4
$ws.Range("N107").Value = -6779.4
$ws.Range("H122").Value = 66675656
$ws.Range("I122").Value = 142867000
$ws.Range("J122").Value = 8224.75
$ws.Range("K122").Value = 428601000
$ws.Range("L122").Value = 24674.25
$ws.Range("M122").Value = -428598550
$ws.Range("N122").Value = -29574.25
$ws.Range("H132").Value = 27905.54
$ws.Range("I132").Value = 2032.9642
$ws.Range("J132").Value = 93763
$ws.Range("K132").Value = 6098.892599999999
$ws.Range("L132").Value = 281289
$ws.Range("M132").Value = -3568.892599999999
$ws.Range("N132").Value = -286349
